$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.832.74'
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').Value = '2.349.47'
$ws.Range('E3').Value = '  -0.92%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '''544.59'
$ws.Range('E5').Value = '  +0.22%  '
$ws.Range('D6').Value = '''136.80'
$ws.Range('E6').Value = '  -2.91%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  -5.04%  '
$ws.Range('D9').Value = '2.347.02'
$ws.Range('E9').Value = '  -0.86%  '
$ws.Range('E10').Value = '  -0.12%  '
$ws.Range('E12').Value = '  -0.59%  '
$ws.Range('D13').Value = '''0.342'
$ws.Range('E13').Value = '  -0.09%  '
$ws.Range('E14').Value = '  -2.67%  '
$ws.Range('D15').Value = '2.774.25'
$ws.Range('E15').Value = '  -0.97%  '
$ws.Range('D16').Value = '60.723.14'
$ws.Range('E16').Value = '  -0.17%  '
$ws.Range('E17').Value = '  -2.05%  '
$ws.Range('D18').Value = '2.345.06'
$ws.Range('E18').Value = '  -1.09%  '
$ws.Range('D19').Value = '''10.64'
$ws.Range('E19').Value = '  +0.82%  '
$ws.Range('E20').Value = '  +0.46%  '
$ws.Range('D21').Value = '''318.83'
$ws.Range('E21').Value = '  +0.54%  '
$ws.Range('D22').Value = '''6.55'
$ws.Range('E22').Value = '  -2.05%  '
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('D24').Value = '''63.35'
$ws.Range('E24').Value = '  +0.73%  '
$ws.Range('E25').Value = '  -6.22%  '
$ws.Range('D26').Value = '''8.32'
$ws.Range('E26').Value = '  +7.68%  '
$ws.Range('D27').Value = '''1.00'
$ws.Range('E27').Value = '  +0.36%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').Value = '''7.96'
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('B29').Value = 'Bittensor'
$ws.Range('C29').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D29').Value = '''499.77'
$ws.Range('E29').Value = '  -3.30%  '
$ws.Range('B30').Value = 'Fetch.AI'
$ws.Range('C30').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D30').Value = '''1.37'
$ws.Range('E30').Value = '  -3.95%  '
$ws.Range('B31').Value = 'PEPE'
$ws.Range('C31').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D31').Value = '0.0₃0859'
$ws.Range('E31').Value = '  -7.00%  '
$ws.Range('B32').Value = 'Kaspa'
$ws.Range('C32').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D32').Value = '''0.145'
$ws.Range('E32').Value = '  +0.45%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').Value = '''1.79'
$ws.Range('E33').Value = '  -2.10%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = '''1.49'
$ws.Range('E34').Value = '  -3.86%  '
$ws.Range('B35').Value = 'FirstDigitalUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D35').Value = '''0.999'
$ws.Range('E35').Value = '  -0.11%  '
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').Value = '''4.59'
$ws.Range('E36').Value = '  -0.98%  '
$ws.Range('B37').Value = 'PolygonEcosystemToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D37').Value = '''0.376'
$ws.Range('E37').Value = '  +0.39%  '
$ws.Range('B38').Value = 'EthereumClassic'
$ws.Range('C38').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D38').Value = '''18.46'
$ws.Range('E38').Value = '  +2.32%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').Value = '''1.83'
$ws.Range('E39').Value = '  +6.28%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D40').Value = '''5.24'
$ws.Range('E40').Value = '  -3.73%  '
$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').Value = '''143.34'
$ws.Range('E41').Value = '  +4.48%  '
$ws.Range('B42').Value = 'USDe'
$ws.Range('C42').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D42').Value = '''0.999'
$ws.Range('E42').Value = '  -0.07%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').Value = '''40.60'
$ws.Range('E43').Value = '  +1.01%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Value = '''143.31'
$ws.Range('E44').Value = '  +3.03%  '
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').Value = '''3.56'
$ws.Range('E45').Value = '  +0.78%  '
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').Value = '''2.04'
$ws.Range('E46').Value = '  -8.75%  '
$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D47').Value = '''0.0518'
$ws.Range('E47').Value = '  +1.02%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').Value = '''19.07'
$ws.Range('E48').Value = '  -6.41%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').Value = '''0.568'
$ws.Range('E49').Value = '  -1.17%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').Value = '''0.0900'
$ws.Range('E50').Value = '  -1.66%  '
$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').Value = '''0.0221'
$ws.Range('E51').Value = '  -1.56%  '
